$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Pediatrics" to "Session"
$ws.Name = "Session"

# The last log entry (row 5, student 201456) was removed from the sheet,
# so delete that entire row and shift the remaining rows up.
$ws.Rows.Item(5).Delete()
